$wb = $excel.ActiveWorkbook

# "Overview" sheet: G2 = Latest HO Xliff Generate Date for first file row.
# This mirrors the de-de hand-off datetime for the same row (both were the
# same shared string in the source workbook), so update it alongside de-de!H2.
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Cells.Item(2, 7).Value = "2016-12-16 09:16:35"

# "zh-cn" sheet: H2 = Correspond Handoff Datetime, L2 = Correspond Handback DateTime
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Cells.Item(2, 8).Value = "2016-12-16 09:16:21"
$wsZhCn.Cells.Item(2, 12).Value = "2016-12-16 09:17:15"

# "de-de" sheet: H2 = Correspond Handoff Datetime (tied to Overview!G2), L2 = Correspond Handback DateTime
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Cells.Item(2, 8).Value = "2016-12-16 09:16:35"
$wsDeDe.Cells.Item(2, 12).Value = "2016-12-16 09:17:34"
